# Update "想去人数" (interest count, column F) values across the four
# sheets to reflect newly scraped numbers, per commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 19
$ws1.Range("F3").Value = 988
$ws1.Range("F4").Value = 230
$ws1.Range("F6").Value = 1140
$ws1.Range("F7").Value = 918
$ws1.Range("F12").Value = 313
$ws1.Range("F17").Value = 1260
$ws1.Range("F18").Value = 2929
$ws1.Range("F19").Value = 153
$ws1.Range("F20").Value = 1541
$ws1.Range("F21").Value = 1298
$ws1.Range("F23").Value = 212
$ws1.Range("F24").Value = 1299
$ws1.Range("F26").Value = 1062
$ws1.Range("F28").Value = 3278
$ws1.Range("F29").Value = 641
$ws1.Range("F31").Value = 1459

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 64
$ws2.Range("F6").Value = 43
$ws2.Range("F8").Value = 10

# Sheet 3: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 772

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 19
$ws4.Range("F3").Value = 772
$ws4.Range("F6").Value = 988
$ws4.Range("F7").Value = 230
$ws4.Range("F10").Value = 1140
$ws4.Range("F11").Value = 918
$ws4.Range("F13").Value = 64
$ws4.Range("F15").Value = 43
$ws4.Range("F16").Value = 43
$ws4.Range("F18").Value = 10
$ws4.Range("F23").Value = 313
$ws4.Range("F28").Value = 1260
$ws4.Range("F29").Value = 2929
$ws4.Range("F30").Value = 153
$ws4.Range("F31").Value = 1541
$ws4.Range("F32").Value = 1298
$ws4.Range("F34").Value = 212
$ws4.Range("F35").Value = 1299
$ws4.Range("F39").Value = 1062
$ws4.Range("F41").Value = 3278
$ws4.Range("F42").Value = 641
$ws4.Range("F44").Value = 1459
